$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C previously held numeric codes mapping each stage to an integer.
# It now holds the same string values as column B (the "mapsto" names),
# with row 5 (stage3) renamed to "sws" and row 6 (REM) renamed to "rem".
$ws.Range("C2").Value = "wake"
$ws.Range("C3").Value = "stage1"
$ws.Range("C4").Value = "stage2"
$ws.Range("C5").Value = "sws"
$ws.Range("C6").Value = "rem"
$ws.Range("C7").Value = "unknown"
